# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update narrative text with the new conversion figures ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.84 = 19190.11 pesos`n✅ 19190.11 pesos = 4.81 = 933.77 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the N/O rate values in rows 10 and 12 ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 206.744
$ws2.Range("O10").Value = 3967.44

$ws2.Range("N12").Value = 3988.99
$ws2.Range("O12").Value = 194.1
